$wb = $excel.ActiveWorkbook

# --- Sheet "신승민" (Sheet1): update task log with new row + revised notes ---
$ws = $wb.Worksheets.Item("신승민")

# Row 2: problem note revised (text changed)
$ws.Cells.Item(2,6).Value = "디비 미완성으로 인해 알고리즘 수정 x // 웹 UI제작"

# Row 3: result note gets a second bullet point appended
$ws.Cells.Item(3,5).Value = "1. detailphone.jsp 구현-> 휴대폰 상세 보기 기능( 이기능에 휴대폰의 스펙이나 성능을 설명해주는 동영상이 있으면 괜찮을것 같아서 youtube링크를 첨부하는쪽으로 추진)  `n2. recomplan.jsp 구현"
$ws.Rows.Item(3).RowHeight = 115.2

# Row 4: brand-new task entry
$ws.Cells.Item(4,1).Value = "web server 구현`n알고리즘 수정"
$ws.Cells.Item(4,2).Value = "웹페이지 구축"
$ws.Cells.Item(4,3).Value = 43606
$ws.Cells.Item(4,3).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(4,4).Value = 43610
$ws.Cells.Item(4,4).NumberFormat = "yyyy\-mm\-dd;@"
$ws.Cells.Item(4,5).Value = "1. 구현했던 jsp파일에 디비에서 데이터를 가져와  웹에 나타내는 코드 구현`n2. 휴대폰 추천 알고리즘 수정(Model.java, recomphone.java 구현)"
$ws.Cells.Item(4,6).Value = "기존코드는 프로그램 시작시 txt파일에서 데이터를 전부다 불러오고 알고리즘을 수행하는 형식이었고 예상변경계획은 디비에서 해당하는 데이터만 가져오는식으로 하려고 했으나 각 휴대폰 기종에 가중치를 주는 형식이기 때문에 부적절하다고 생각되서 요즘제 추천을 받을때 디비에서 휴대폰을 전부 가져오고 우선순위를 계산하는 형식으로 코드 변경"
$ws.Rows.Item(4).RowHeight = 249.6

# Give column A a bit more breathing room now that it holds longer task names
$ws.Columns.Item(1).ColumnWidth = 16.7

# Move the active selection
$ws.Range("K3").Select()

# --- Sheet "박성곤" (Sheet3): clear stale explicit row heights on empty rows ---
$ws3 = $wb.Worksheets.Item("박성곤")
$ws3.Rows.Item(1).EntireRow.AutoFit()
$ws3.Range("A4:F13").EntireRow.AutoFit()

$wb.Save()
